$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KODE_REKSADANA cell (N2) from RD00014 to RD00015
$ws.Range("N2").Value = "RD00015"

# Update PREPARATION cell (F2) text to reflect the new Kode Reksadana value
$ws.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nKode Reksadana : RD00015"

# Update the active selection from O2 to G2
$ws.Range("G2").Select()
